# Refresh the cryptocurrency Price (D) and Volume(1h) (E) columns for
# rows 2-51 with the latest scraped figures (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> column -> new text value. Price strings that look like a
# plain decimal number (e.g. "1.005") are pre-escaped here with a leading
# apostrophe so Excel stores them as text/quote-prefixed, exactly like the
# source data, instead of silently re-typing them as numbers.
$updates = @{
    2 = @{ "D" = '29.924.95'; "E" = '  +1.76%  ' }
    3 = @{ "D" = '1.939.18'; "E" = '  +1.35%  ' }
    4 = @{ "E" = '  -0.42%  ' }
    5 = @{ "E" = '  +2.99%  ' }
    6 = @{ "D" = '''1.005'; "E" = '  -0.37%  ' }
    7 = @{ "D" = '''0.4839'; "E" = '  +0.49%  ' }
    8 = @{ "D" = '''0.4127'; "E" = '  +1.28%  ' }
    9 = @{ "D" = '''0.08178'; "E" = '  -0.39%  ' }
    10 = @{ "D" = '''1.016'; "E" = '  -0.39%  ' }
    11 = @{ "D" = '''23.66'; "E" = '  +0.92%  ' }
    12 = @{ "D" = '1.957.45'; "E" = '  +1.97%  ' }
    13 = @{ "D" = '''6.094'; "E" = '  +0.97%  ' }
    14 = @{ "D" = '''7.303'; "E" = '  +1.19%  ' }
    15 = @{ "D" = '''91.19'; "E" = '  +0.01%  ' }
    16 = @{ "D" = '''0.06855'; "E" = '  +0.97%  ' }
    17 = @{ "D" = '''1.007'; "E" = '  -0.25%  ' }
    18 = @{ "D" = '''0.00001037'; "E" = '  -0.16%  ' }
    19 = @{ "D" = '''17.80'; "E" = '  +0.27%  ' }
    20 = @{ "E" = '  -0.22%  ' }
    21 = @{ "D" = '29.917.44'; "E" = '  +1.54%  ' }
    22 = @{ "D" = '''5.640'; "E" = '  +0.06%  ' }
    23 = @{ "D" = '''11.89'; "E" = '  +1.32%  ' }
    24 = @{ "E" = '  -0.33%  ' }
    25 = @{ "D" = '2.185.45'; "E" = '  +1.49%  ' }
    26 = @{ "D" = '''6.699'; "E" = '  +1.06%  ' }
    27 = @{ "D" = '''156.62'; "E" = '  -0.13%  ' }
    28 = @{ "D" = '''20.09'; "E" = '  +0.41%  ' }
    29 = @{ "D" = '''2.104'; "E" = '  -0.17%  ' }
    30 = @{ "D" = '''121.40'; "E" = '  +1.05%  ' }
    31 = @{ "D" = '''1.010'; "E" = '  -1.21%  ' }
    32 = @{ "D" = '''0.09634'; "E" = '  +0.83%  ' }
    33 = @{ "D" = '''5.594'; "E" = '  +1.18%  ' }
    34 = @{ "D" = '''1.426'; "E" = '  +3.72%  ' }
    35 = @{ "D" = '''3.545'; "E" = '  -0.42%  ' }
    36 = @{ "D" = '''0.06569'; "E" = '  +7.08%  ' }
    37 = @{ "D" = '''0.02288'; "E" = '  +0.11%  ' }
    38 = @{ "D" = '''1.211'; "E" = '  +2.64%  ' }
    39 = @{ "D" = '''0.5974'; "E" = '  +0.05%  ' }
    40 = @{ "D" = '''8.003'; "E" = '  -0.49%  ' }
    41 = @{ "E" = '  -0.37%  ' }
    42 = @{ "D" = '''0.1852'; "E" = '  +0.11%  ' }
    43 = @{ "D" = '''2.498'; "E" = '  +4.02%  ' }
    44 = @{ "D" = '''1.271'; "E" = '  -0.72%  ' }
    45 = @{ "D" = '''12.40'; "E" = '  -0.15%  ' }
    46 = @{ "D" = '''0.07495'; "E" = '  -1.33%  ' }
    47 = @{ "D" = '''0.5574'; "E" = '  +0.14%  ' }
    48 = @{ "D" = '''1.991'; "E" = '  +1.85%  ' }
    49 = @{ "D" = '''117.35'; "E" = '  +0.14%  ' }
    50 = @{ "D" = '''72.70'; "E" = '  +0.14%  ' }
    51 = @{ "D" = '''2.415'; "E" = '  -0.62%  ' }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $cellRef = "$col$row"
        $ws.Range($cellRef).Value = $updates[$row][$col]
    }
}
